$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '37.543.59'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = "'" + '2.030.26'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'" + '255.47'
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("D7").Value = "'" + '1.00'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'" + '57.52'
$ws.Range("E8").Value = '  -5.86%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").Value = "'" + '0.101'
$ws.Range("E11").Value = '  -2.05%  '
$ws.Range("D12").Value = "'" + '14.61'
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").Value = "'" + '2.330.09'
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").Value = "'" + '0.818'
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = "'" + '21.30'
$ws.Range("E15").Value = '  -4.21%  '
$ws.Range("D16").Value = "'" + '5.37'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").Value = "'" + '2.017.95'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = "'" + '37.468.14'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = "'" + '69.71'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("E21").Value = '  +1.22%  '
$ws.Range("D22").Value = "'" + '229.18'
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = "'" + '2.62'
$ws.Range("E23").Value = '  +3.98%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = "'" + '1.00'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").Value = "'" + '164.54'
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").Value = "'" + '9.09'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = "'" + '19.96'
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("E29").Value = '  -13.28%  '
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").Value = "'" + '0.0666'
$ws.Range("E32").Value = '  +7.16%  '
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("E35").Value = '  +6.44%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("D38").Value = "'" + '3.39'
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").Value = "'" + '5.35'
$ws.Range("E39").Value = '  -2.43%  '
$ws.Range("E40").Value = '  +3.13%  '
$ws.Range("D41").Value = "'" + '0.0966'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = "'" + '1.401.52'
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = "'" + '16.09'
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").Value = "'" + '91.35'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").Value = "'" + '7.37'
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("D50").Value = "'" + '2.02'
$ws.Range("E50").Value = '  +3.01%  '
$ws.Range("D51").Value = "'" + '2.221.75'
$ws.Range("E51").Value = '  +2.61%  '
